$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F13").Value = "verger"
$ws.Range("F13").Select()
